$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "Manufacturer Part Number" -> "Manufacturer Part #" ---
$ws.Range("I4").Value = "Manufacturer Part #"

# --- Row 7: merge the two power-cord rows into a single row ---
# Update the primary cord row text (drop the "Canada, USA" suffix)
$ws.Range("C7").Value = "CORD 18AWG NEMA 1-15P - C7 6.56'"
# Store the Digi-Key catalogue number as a real number now
$ws.Range("I7").Value = 6010.5274
# Corrected supplier part number
$ws.Range("K7").Value = "485-4157-ND"

# Remove the now-redundant "alternative" cord row (row 8) entirely
$ws.Rows(8).Delete()

# --- Former row 10 ("FUSE 1A 250VAC 5x20"), now row 9 after the delete ---
$ws.Range("C9").Value = "FUSE"

# --- Former row 11 (fuse-holder cover), now row 10 ---
$ws.Range("K10").Value = "35-4527C-ND"

# --- Former row 12 (QC receptacle connector), now row 11 ---
$ws.Range("C11").Value = "CONN QC RCPT 16-20 AWG"

# --- Remove the old blank placeholder row (former row 13, now row 12) ---
$ws.Rows(12).Delete()

# --- Add the new "Washer" line item as row 12 ---
$ws.Range("A12").Value = 4
$ws.Range("B12").Value = "none"
$ws.Range("C12").Value = "WASHER FLAT RETAINING #4 NYLON"
$ws.Range("D12").Value = "N.A."
$ws.Range("H12").Value = "Essentra Components"
$ws.Range("I12").Value = "16FWRT004050"
$ws.Range("J12").Value = "Digi-Key"
$ws.Range("K12").Value = "RPC6339-ND"

# Row 5 picked up a taller, manually-set row height
$ws.Rows(5).RowHeight = 75

# Column width tweaks left behind by the edit
$ws.Columns(1).AutoFit()
$ws.Columns(3).ColumnWidth = 36.7109375
$ws.Columns(8).ColumnWidth = 19.7109375
$ws.Columns(9).ColumnWidth = 24
$ws.Columns(10).ColumnWidth = 9.7109375
$ws.Columns(11).ColumnWidth = 21.85546875

# Scroll position / selection as left by the editor
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("C22").Select()
